$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header row) updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 updates
$ws.Range("B2").Value = 18.42673536815974
$ws.Range("C2").Value = 4.737029407806765
$ws.Range("D2").Value = 4.8920507547112022
$ws.Range("E2").Value = 1.2270924346332597

# Row 3 updates
$ws.Range("B3").Value = 32.582619160327134
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = -5.1703741057961867
$ws.Range("E3").Value = 6.6860536061983291

# Restore selection to reflect the updated data range
$ws.Range("B1:E3").Select()
